$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (years 2000-2009), which shifts rows 12-21 (years 2010-2019) up to rows 2-11
$ws.Rows("2:11").Delete()
